# Updates Anxa1-Dysf NATMI TPM output values for all 25 data rows (rows 2-26)
# as described by the commit "update scripts wuth new tpm"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 22.61022533333333
$ws.Range("H2").Value = 67.830676
$ws.Range("I2").Value = 0.04352672200082041
$ws.Range("J2").Value = 0.04795217939334551
$ws.Range("M2").Value = 32.92864966666667
$ws.Range("N2").Value = 98.785949
$ws.Range("O2").Value = 0.8913028757746132
$ws.Range("P2").Value = 0.9086533747458582
$ws.Range("Q2").Value = 744.5241888857248
$ws.Range("R2").Value = 6700.717699971524
$ws.Range("S2").Value = 0.03879549249237336
$ws.Range("T2").Value = 0.04357190963218219

# Row 3
$ws.Range("G3").Value = 22.61022533333333
$ws.Range("H3").Value = 67.830676
$ws.Range("I3").Value = 0.04352672200082041
$ws.Range("J3").Value = 0.04795217939334551
$ws.Range("O3").Value = 0.04855270847140305
$ws.Range("P3").Value = 0.04949785713105716
$ws.Range("Q3").Value = 40.55710676515022
$ws.Range("R3").Value = 365.013960886352
$ws.Range("S3").Value = 0.002113340244021638
$ws.Range("T3").Value = 0.002373530124734639

# Row 4
$ws.Range("G4").Value = 22.61022533333333
$ws.Range("H4").Value = 67.830676
$ws.Range("I4").Value = 0.04352672200082041
$ws.Range("J4").Value = 0.04795217939334551
$ws.Range("M4").Value = 0.05066733333333334
$ws.Range("N4").Value = 0.152002
$ws.Range("O4").Value = 0.001371448278777914
$ws.Range("P4").Value = 0.001398145502131279
$ws.Range("Q4").Value = 1.145599823705778
$ws.Range("R4").Value = 10.310398413352
$ws.Range("S4").Value = 0.00005969464796886991
$ws.Range("T4").Value = 0.00006704412393619824

# Row 5
$ws.Range("G5").Value = 22.61022533333333
$ws.Range("H5").Value = 67.830676
$ws.Range("I5").Value = 0.04352672200082041
$ws.Range("J5").Value = 0.04795217939334551
$ws.Range("M5").Value = 2.116331
$ws.Range("N5").Value = 4.232661999999999
$ws.Range("O5").Value = 0.05728421679861464
$ws.Range("P5").Value = 0.03893289126025962
$ws.Range("Q5").Value = 47.85072078991866
$ws.Range("R5").Value = 287.1043247395119
$ws.Range("S5").Value = 0.002493394179628026
$ws.Range("T5").Value = 0.001866916986013583

# Row 6
$ws.Range("G6").Value = 22.61022533333333
$ws.Range("H6").Value = 67.830676
$ws.Range("I6").Value = 0.04352672200082041
$ws.Range("J6").Value = 0.04795217939334551
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.055001
$ws.Range("N6").Value = 0.165003
$ws.Range("O6").Value = 0.001488750676591046
$ws.Range("P6").Value = 0.001517731360693724
$ws.Range("Q6").Value = 1.243585003558667
$ws.Range("R6").Value = 11.192265032028
$ws.Range("S6").Value = 0.00006480043682851173
$ws.Range("T6").Value = 0.00007277852647889184

# Row 7
$ws.Range("G7").Value = 178.3379163333334
$ws.Range("H7").Value = 535.0137490000001
$ws.Range("I7").Value = 0.3433165654922813
$ws.Range("J7").Value = 0.3782223144872436
$ws.Range("M7").Value = 32.92864966666667
$ws.Range("N7").Value = 98.785949
$ws.Range("O7").Value = 0.8913028757746132
$ws.Range("P7").Value = 0.9086533747458582
$ws.Range("Q7").Value = 5872.426769223645
$ws.Range("R7").Value = 52851.84092301281
$ws.Range("S7").Value = 0.3059990421243337
$ws.Range("T7").Value = 0.3436729824630232

# Row 8
$ws.Range("G8").Value = 178.3379163333334
$ws.Range("H8").Value = 535.0137490000001
$ws.Range("I8").Value = 0.3433165654922813
$ws.Range("J8").Value = 0.3782223144872436
$ws.Range("O8").Value = 0.04855270847140305
$ws.Range("P8").Value = 0.04949785713105716
$ws.Range("Q8").Value = 319.893756314861
$ws.Range("R8").Value = 2879.043806833748
$ws.Range("S8").Value = 0.01666894911775009
$ws.Range("T8").Value = 0.01872119408626735

# Row 9
$ws.Range("G9").Value = 178.3379163333334
$ws.Range("H9").Value = 535.0137490000001
$ws.Range("I9").Value = 0.3433165654922813
$ws.Range("J9").Value = 0.3782223144872436
$ws.Range("M9").Value = 0.05066733333333334
$ws.Range("N9").Value = 0.152002
$ws.Range("O9").Value = 0.001371448278777914
$ws.Range("P9").Value = 0.001398145502131279
$ws.Range("Q9").Value = 9.035906652833113
$ws.Range("R9").Value = 81.323159875498
$ws.Range("S9").Value = 0.0004708409128203342
$ws.Range("T9").Value = 0.0005288098278060218

# Row 10
$ws.Range("G10").Value = 178.3379163333334
$ws.Range("H10").Value = 535.0137490000001
$ws.Range("I10").Value = 0.3433165654922813
$ws.Range("J10").Value = 0.3782223144872436
$ws.Range("M10").Value = 2.116331
$ws.Range("N10").Value = 4.232661999999999
$ws.Range("O10").Value = 0.05728421679861464
$ws.Range("P10").Value = 0.03893289126025962
$ws.Range("Q10").Value = 377.4220608116397
$ws.Range("R10").Value = 2264.532364869838
$ws.Range("S10").Value = 0.01966662056821563
$ws.Range("T10").Value = 0.01472528824213558

# Row 11
$ws.Range("G11").Value = 178.3379163333334
$ws.Range("H11").Value = 535.0137490000001
$ws.Range("I11").Value = 0.3433165654922813
$ws.Range("J11").Value = 0.3782223144872436
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.055001
$ws.Range("N11").Value = 0.165003
$ws.Range("O11").Value = 0.001488750676591046
$ws.Range("P11").Value = 0.001517731360693724
$ws.Range("Q11").Value = 9.808763736249668
$ws.Range("R11").Value = 88.27887362624702
$ws.Range("S11").Value = 0.0005111127691615479
$ws.Range("T11").Value = 0.000574039868011454

# Row 12
$ws.Range("G12").Value = 90.63663000000001
$ws.Range("H12").Value = 271.90989
$ws.Range("I12").Value = 0.174483683330882
$ws.Range("J12").Value = 0.1922238225092264
$ws.Range("M12").Value = 32.92864966666667
$ws.Range("N12").Value = 98.785949
$ws.Range("O12").Value = 0.8913028757746132
$ws.Range("P12").Value = 0.9086533747458582
$ws.Range("Q12").Value = 2984.54183623729
$ws.Range("R12").Value = 26860.87652613561
$ws.Range("S12").Value = 0.1555178087285621
$ws.Range("T12").Value = 0.1746648250295574

# Row 13
$ws.Range("G13").Value = 90.63663000000001
$ws.Range("H13").Value = 271.90989
$ws.Range("I13").Value = 0.174483683330882
$ws.Range("J13").Value = 0.1922238225092264
$ws.Range("O13").Value = 0.04855270847140305
$ws.Range("P13").Value = 0.04949785713105716
$ws.Range("Q13").Value = 162.57951548692
$ws.Range("R13").Value = 1463.21563938228
$ws.Range("S13").Value = 0.008471655409780923
$ws.Range("T13").Value = 0.009514667303747379

# Row 14
$ws.Range("G14").Value = 90.63663000000001
$ws.Range("H14").Value = 271.90989
$ws.Range("I14").Value = 0.174483683330882
$ws.Range("J14").Value = 0.1922238225092264
$ws.Range("M14").Value = 0.05066733333333334
$ws.Range("N14").Value = 0.152002
$ws.Range("O14").Value = 0.001371448278777914
$ws.Range("P14").Value = 0.001398145502131279
$ws.Range("Q14").Value = 4.59231634442
$ws.Range("R14").Value = 41.33084709978
$ws.Range("S14").Value = 0.0002392953471789687
$ws.Range("T14").Value = 0.0002687568728437563

# Row 15
$ws.Range("G15").Value = 90.63663000000001
$ws.Range("H15").Value = 271.90989
$ws.Range("I15").Value = 0.174483683330882
$ws.Range("J15").Value = 0.1922238225092264
$ws.Range("M15").Value = 2.116331
$ws.Range("N15").Value = 4.232661999999999
$ws.Range("O15").Value = 0.05728421679861464
$ws.Range("P15").Value = 0.03893289126025962
$ws.Range("Q15").Value = 191.81710980453
$ws.Range("R15").Value = 1150.90265882718
$ws.Range("S15").Value = 0.009995161143747069
$ws.Range("T15").Value = 0.007483829179383159

# Row 16
$ws.Range("G16").Value = 90.63663000000001
$ws.Range("H16").Value = 271.90989
$ws.Range("I16").Value = 0.174483683330882
$ws.Range("J16").Value = 0.1922238225092264
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.055001
$ws.Range("N16").Value = 0.165003
$ws.Range("O16").Value = 0.001488750676591046
$ws.Range("P16").Value = 0.001517731360693724
$ws.Range("Q16").Value = 4.985105286630001
$ws.Range("R16").Value = 44.86594757967001
$ws.Range("S16").Value = 0.0002597627016129484
$ws.Range("T16").Value = 0.0002917441236946771

# Row 17
$ws.Range("G17").Value = 143.820236
$ws.Range("H17").Value = 287.640472
$ws.Range("I17").Value = 0.2768669192002915
$ws.Range("J17").Value = 0.2033443911738485
$ws.Range("M17").Value = 32.92864966666667
$ws.Range("N17").Value = 98.785949
$ws.Range("O17").Value = 0.8913028757746132
$ws.Range("P17").Value = 0.9086533747458582
$ws.Range("Q17").Value = 4735.806166221321
$ws.Range("R17").Value = 28414.83699732793
$ws.Range("S17").Value = 0.2467722812900773
$ws.Range("T17").Value = 0.1847695672757593

# Row 18
$ws.Range("G18").Value = 143.820236
$ws.Range("H18").Value = 287.640472
$ws.Range("I18").Value = 0.2768669192002915
$ws.Range("J18").Value = 0.2033443911738485
$ws.Range("O18").Value = 0.04855270847140305
$ws.Range("P18").Value = 0.04949785713105716
$ws.Range("Q18").Value = 257.9776442051573
$ws.Range("R18").Value = 1547.865865230944
$ws.Range("S18").Value = 0.01344263881330726
$ws.Range("T18").Value = 0.01006511162272495

# Row 19
$ws.Range("G19").Value = 143.820236
$ws.Range("H19").Value = 287.640472
$ws.Range("I19").Value = 0.2768669192002915
$ws.Range("J19").Value = 0.2033443911738485
$ws.Range("M19").Value = 0.05066733333333334
$ws.Range("N19").Value = 0.152002
$ws.Range("O19").Value = 0.001371448278777914
$ws.Range("P19").Value = 0.001398145502131279
$ws.Range("Q19").Value = 7.286987837490667
$ws.Range("R19").Value = 43.72192702494399
$ws.Range("S19").Value = 0.0003797086597877835
$ws.Range("T19").Value = 0.0002843050459033397

# Row 20
$ws.Range("G20").Value = 143.820236
$ws.Range("H20").Value = 287.640472
$ws.Range("I20").Value = 0.2768669192002915
$ws.Range("J20").Value = 0.2033443911738485
$ws.Range("M20").Value = 2.116331
$ws.Range("N20").Value = 4.232661999999999
$ws.Range("O20").Value = 0.05728421679861464
$ws.Range("P20").Value = 0.03893289126025962
$ws.Range("Q20").Value = 304.371223874116
$ws.Range("R20").Value = 1217.484895496464
$ws.Range("S20").Value = 0.01586010462383402
$ws.Range("T20").Value = 0.00791678506995514

# Row 21
$ws.Range("G21").Value = 143.820236
$ws.Range("H21").Value = 287.640472
$ws.Range("I21").Value = 0.2768669192002915
$ws.Range("J21").Value = 0.2033443911738485
$ws.Range("K21").Value = 1
$ws.Range("L21").Value = 0.3333333333333333
$ws.Range("M21").Value = 0.055001
$ws.Range("N21").Value = 0.165003
$ws.Range("O21").Value = 0.001488750676591046
$ws.Range("P21").Value = 0.001517731360693724
$ws.Range("Q21").Value = 7.910256800236
$ws.Range("R21").Value = 47.461540801416
$ws.Range("S21").Value = 0.0004121858132851123
$ws.Range("T21").Value = 0.000308622159505722

# Row 22
$ws.Range("G22").Value = 84.051186
$ws.Range("H22").Value = 252.153558
$ws.Range("I22").Value = 0.1618061099757246
$ws.Range("J22").Value = 0.1782572924363359
$ws.Range("M22").Value = 32.92864966666667
$ws.Range("N22").Value = 98.785949
$ws.Range("O22").Value = 0.8913028757746132
$ws.Range("P22").Value = 0.9086533747458582
$ws.Range("Q22").Value = 2767.692057861838
$ws.Range("R22").Value = 24909.22852075654
$ws.Range("S22").Value = 0.1442182511392667
$ws.Range("T22").Value = 0.1619740903453359

# Row 23
$ws.Range("G23").Value = 84.051186
$ws.Range("H23").Value = 252.153558
$ws.Range("I23").Value = 0.1618061099757246
$ws.Range("J23").Value = 0.1782572924363359
$ws.Range("O23").Value = 0.04855270847140305
$ws.Range("P23").Value = 0.04949785713105716
$ws.Range("Q23").Value = 150.766870921624
$ws.Range("R23").Value = 1356.901838294616
$ws.Range("S23").Value = 0.00785612488654314
$ws.Range("T23").Value = 0.00882335399358283

# Row 24
$ws.Range("G24").Value = 84.051186
$ws.Range("H24").Value = 252.153558
$ws.Range("I24").Value = 0.1618061099757246
$ws.Range("J24").Value = 0.1782572924363359
$ws.Range("M24").Value = 0.05066733333333334
$ws.Range("N24").Value = 0.152002
$ws.Range("O24").Value = 0.001371448278777914
$ws.Range("P24").Value = 0.001398145502131279
$ws.Range("Q24").Value = 4.258649458124
$ws.Range("R24").Value = 38.327845123116
$ws.Range("S24").Value = 0.0002219087110219574
$ws.Range("T24").Value = 0.0002492296316419631

# Row 25
$ws.Range("G25").Value = 84.051186
$ws.Range("H25").Value = 252.153558
$ws.Range("I25").Value = 0.1618061099757246
$ws.Range("J25").Value = 0.1782572924363359
$ws.Range("M25").Value = 2.116331
$ws.Range("N25").Value = 4.232661999999999
$ws.Range("O25").Value = 0.05728421679861464
$ws.Range("P25").Value = 0.03893289126025962
$ws.Range("Q25").Value = 177.880130518566
$ws.Range("R25").Value = 1067.280783111396
$ws.Range("S25").Value = 0.009268936283189894
$ws.Range("T25").Value = 0.006940071782772166

# Row 26
$ws.Range("G26").Value = 84.051186
$ws.Range("H26").Value = 252.153558
$ws.Range("I26").Value = 0.1618061099757246
$ws.Range("J26").Value = 0.1782572924363359
$ws.Range("K26").Value = 1
$ws.Range("L26").Value = 0.3333333333333333
$ws.Range("M26").Value = 0.055001
$ws.Range("N26").Value = 0.165003
$ws.Range("O26").Value = 0.001488750676591046
$ws.Range("P26").Value = 0.001517731360693724
$ws.Range("Q26").Value = 4.622899281186
$ws.Range("R26").Value = 41.606093530674
$ws.Range("S26").Value = 0.0002408889557029252
$ws.Range("T26").Value = 0.0002705466830029792
